$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: 12-Aug-2022, 7, Chikkamagaluru
$ws.Range("A6").Value = 44785
$ws.Range("A6").NumberFormat = "d-mmm"
$ws.Range("B6").Value = 7
$ws.Range("C6").Value = "Chikkamagaluru"

# Row 7: 21-Jul-2022, 3, Bengaluru
$ws.Range("A7").Value = 44763
$ws.Range("A3").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "Bengaluru"

# Update the active selection to D7, matching the target workbook state
$ws.Range("D7").Select()
